$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 ("Prescaler 8" row), shifting
# the existing data rows down by one. Excel automatically adjusts the
# relative formulas in the shifted rows (E2->E3, G2->G3, etc.).
$ws.Rows("2:2").Insert()

# Fill in the new row 2 (Prescaler = 1)
$ws.Range("A2").Value = 1
$ws.Range("B2").Formula = "=E3/A2*1000000"
$ws.Range("C2").Formula = "=G3/D2-1"
$ws.Range("D2").Formula = "=1/B2"

# Update the "Overflow(S)" value that now lives on row 3 together with the
# Frequency/Mhz columns (E3/F3/G3)
$ws.Range("G3").Value = 4.1

# Column width tweaks: column D no longer relies on bestFit sizing and
# column G gets an explicit width now that it holds data again.
$ws.Columns("D").ColumnWidth = 10.33
$ws.Columns("G").ColumnWidth = 11.15
